# Add "Program" and "Lot" columns to the end of the "ppv" table on the PPV sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PPV")
$lo = $ws.ListObjects.Item("ppv")

$colProgram = $lo.ListColumns.Add()
$colProgram.Range.Cells.Item(1, 1).Value = "Program"

$colLot = $lo.ListColumns.Add()
$colLot.Range.Cells.Item(1, 1).Value = "Lot"

$ws.Range("J2").Select()
